$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 2801.106610098821

# Row 3
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 37.47995879822157

# Row 4
$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 41249014.21622031
$ws.Range("D4").Value = 3.223369029078222
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 41249034.57576305
